$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 2798
$ws.Range("J3").Value = 2894
$ws.Range("H4").Value = 1695
$ws.Range("J6").Value = 3540
$ws.Range("H7").Value = 26005

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 36
$ws.Range("J3").Value = 27
$ws.Range("J7").Value = 108

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("J2").Value = 10
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J3").Value = 12
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J3").Value = 48
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 146

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 77
$ws.Range("J3").Value = 146
$ws.Range("J6").Value = 112
$ws.Range("J7").Value = 367

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 27
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 78
$ws.Range("J6").Value = 92
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 79
$ws.Range("J5").Value = 27
$ws.Range("J7").Value = 302
$ws.Range("J8").Value = 641
$ws.Range("J11").Value = 145
$ws.Range("J14").Value = 36
$ws.Range("J18").Value = 105
$ws.Range("J19").Value = 320
$ws.Range("J20").Value = 209
$ws.Range("J23").Value = 105
$ws.Range("J26").Value = 14
$ws.Range("J29").Value = 579
$ws.Range("J30").Value = 40
$ws.Range("J33").Value = 419
$ws.Range("J42").Value = 400
$ws.Range("J46").Value = 36
$ws.Range("J49").Value = 63
$ws.Range("J50").Value = 58
$ws.Range("J51").Value = 135
$ws.Range("J52").Value = 259
$ws.Range("J54").Value = 199
$ws.Range("J55").Value = 123
$ws.Range("J56").Value = 11
$ws.Range("H63").Value = 246
$ws.Range("J65").Value = 259
$ws.Range("J66").Value = 23
$ws.Range("J67").Value = 367
$ws.Range("J73").Value = 92
$ws.Range("J76").Value = 145
$ws.Range("J77").Value = 86
$ws.Range("J79").Value = 301
$ws.Range("J83").Value = 238
$ws.Range("J84").Value = 92
$ws.Range("J85").Value = 467
$ws.Range("J86").Value = 59
$ws.Range("J89").Value = 108
$ws.Range("J90").Value = 112
$ws.Range("J95").Value = 161
$ws.Range("J96").Value = 118
$ws.Range("J97").Value = 63
$ws.Range("J99").Value = 146
$ws.Range("H101").Value = 26005

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 70
$ws.Range("J3").Value = 84
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 238

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J2").Value = 62
$ws.Range("J6").Value = 43
$ws.Range("J7").Value = 161

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 110
$ws.Range("J3").Value = 127
$ws.Range("J4").Value = 22
$ws.Range("J6").Value = 144
$ws.Range("J7").Value = 419

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J2").Value = 14
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 199

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 196
$ws.Range("J4").Value = 33
$ws.Range("J7").Value = 579

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 77
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 320

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J4").Value = 14
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 111
$ws.Range("J3").Value = 178
$ws.Range("J6").Value = 133
$ws.Range("J7").Value = 467

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J6").Value = 208
$ws.Range("J7").Value = 400

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J3").Value = 23
$ws.Range("J7").Value = 123

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J2").Value = 13
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J2").Value = 34
$ws.Range("J7").Value = 105

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 34
$ws.Range("J6").Value = 19

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 111
$ws.Range("J7").Value = 301

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 68
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 209

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 30
$ws.Range("J7").Value = 105

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J3").Value = 70
$ws.Range("J6").Value = 115
$ws.Range("J7").Value = 259

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J3").Value = 19
$ws.Range("J7").Value = 58

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 23

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 145

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J6").Value = 24
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 198
$ws.Range("J4").Value = 35
$ws.Range("J6").Value = 186
$ws.Range("J7").Value = 641

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J4").Value = 28
$ws.Range("J7").Value = 59

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 112

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J3").Value = 28
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 11

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 302
